$wb = $excel.ActiveWorkbook

# Fix the "Sonar" tag label spelling across the workbook: סונאר -> סונר
$wsContacts = $wb.Worksheets.Item("contacts_master")
$wsMembership = $wb.Worksheets.Item("tag_membership")
$wsMeta = $wb.Worksheets.Item("tag_meta")

# contacts_master!C5: CPO Amir Example's label changes from מבצעים to סונר
$wsContacts.Range("C5").Value = "סונר"

# tag_membership!A6: סונאר -> סונר
$wsMembership.Range("A6").Value = "סונר"

# tag_meta!A6: סונאר -> סונר
$wsMeta.Range("A6").Value = "סונר"

# Update selections on each sheet to match the saved state
$wsContacts.Range("C6").Select() | Out-Null
$wsMembership.Range("C15").Select() | Out-Null
$wsMeta.Range("B21").Select() | Out-Null

# tag_meta becomes the active sheet/tab
$wsMeta.Activate() | Out-Null
